# The sheet holds yearly employment figures (2000年..2020年) in columns
# A (year label), B (rural employment), C (urban employment), D (total).
# The update drops the 2000-2009 rows, which shifts the 2010-2020 rows
# up by ten positions, and appends two new years: 2021年 and 2022年.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows for 2000年-2009年 (original rows 2-11). Everything below
# (2010年-2020年, originally rows 12-22) shifts up to rows 2-12.
$ws.Rows("2:11").Delete()

# Row 12 now holds 2020年 - the last existing data row. Duplicate its
# formatting (bold font + border, same as every other year cell in
# column A) down into the two new rows before overwriting the values.
$ws.Range("A12").Copy($ws.Range("A13"))
$ws.Range("A12").Copy($ws.Range("A14"))

# New row for 2021年
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 27879
$ws.Range("C13").Value = 46773
$ws.Range("D13").Value = 74652

# New row for 2022年
$ws.Range("A14").Value = "2022年"
$ws.Range("B14").Value = 27420
$ws.Range("C14").Value = 45931
$ws.Range("D14").Value = 73351
